$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the old
# N/O/P ("Late" / "heading" / "Outstanding") columns one place to the
# right (-> O/P/Q). This mirrors a manual "Insert Column" in Excel.
$ws.Columns("N").Insert()

# The inserted column picks up the formatting (including width) of the
# column to its left (M, "In Advance" = 10.7109375 characters wide).
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and move its selection,
# and clear the "Transactions" sheet's previous tab-selected state.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
